$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws1.Range("H113").Value = 13902638
$ws1.Range("J113").Value = 15713.571
$ws1.Range("L113").Value = 15713.571
$ws1.Range("N113").Value = -22221.571
$ws1.Range("H118").Value = 3287
$ws1.Range("I118").Value = 979.6667
$ws1.Range("K118").Value = 2939.0001
$ws1.Range("M118").Value = -1282.0001
$ws1.Range("H138").Value = 5179.0967
$ws1.Range("I138").Value = 1823.5
$ws1.Range("J138").Value = 6346.2607
$ws1.Range("K138").Value = 5470.5
$ws1.Range("L138").Value = 19038.7821
$ws1.Range("M138").Value = -330.5
$ws1.Range("N138").Value = -29318.7821

$ws2 = $wb.Worksheets.Item("ARM")
$ws2.Range("H63").Value = 1866.3334
$ws2.Range("I63").Value = 1866.3334
$ws2.Range("J63").Value = 0
$ws2.Range("K63").Value = 1866.3334
$ws2.Range("L63").Value = 0
$ws2.Range("M63").Value = -1180.3334
$ws2.Range("N63").ClearContents()
$ws2.Range("H66").Value = 1866.3334
$ws2.Range("I66").Value = 1866.3334
$ws2.Range("J66").Value = 0
$ws2.Range("K66").Value = 9331.666999999999
$ws2.Range("L66").Value = 0
$ws2.Range("M66").Value = -5899.666999999999
$ws2.Range("N66").ClearContents()
$ws2.Range("H110").Value = 1310
$ws2.Range("I110").Value = 1270
$ws2.Range("J110").Value = 1450
$ws2.Range("K110").Value = 1270
$ws2.Range("L110").Value = 1450
$ws2.Range("M110").Value = 775
$ws2.Range("N110").Value = -5540
$ws2.Range("H132").Value = 4858.6724
$ws2.Range("I132").Value = 3492.3333
$ws2.Range("J132").Value = 8701.5
$ws2.Range("K132").Value = 10476.9999
$ws2.Range("L132").Value = 26104.5
$ws2.Range("M132").Value = -7946.999899999999
$ws2.Range("N132").Value = -31164.5

$ws3 = $wb.Worksheets.Item("BSM")
$ws3.Range("H81").Value = 66397.60000000001
$ws3.Range("J81").Value = 66397.60000000001
$ws3.Range("L81").Value = 66397.60000000001
$ws3.Range("N81").Value = -68519.60000000001
$ws3.Range("H84").Value = 66397.60000000001
$ws3.Range("J84").Value = 66397.60000000001
$ws3.Range("L84").Value = 199192.8
$ws3.Range("N84").Value = -209800.8
$ws3.Range("H86").Value = 73591740
$ws3.Range("I86").Value = 22819940
$ws3.Range("J86").Value = 166673390
$ws3.Range("K86").Value = 22819940
$ws3.Range("L86").Value = 166673390
$ws3.Range("M86").Value = -22818817
$ws3.Range("N86").Value = -166675636
$ws3.Range("H89").Value = 73591740
$ws3.Range("I89").Value = 22819940
$ws3.Range("J89").Value = 166673390
$ws3.Range("K89").Value = 114099700
$ws3.Range("L89").Value = 833366950
$ws3.Range("M89").Value = -114094084
$ws3.Range("N89").Value = -833378182
$ws3.Range("H94").Value = 4596
$ws3.Range("I94").Value = 1759.3334
$ws3.Range("J94").Value = 8000
$ws3.Range("K94").Value = 1759.3334
$ws3.Range("L94").Value = 8000
$ws3.Range("M94").Value = -1308.3334
$ws3.Range("N94").Value = -8902

$ws4 = $wb.Worksheets.Item("CRP")
$ws4.Range("H99").Value = 9941.666999999999
$ws4.Range("I99").Value = 11216.667
$ws4.Range("K99").Value = 11216.667
$ws4.Range("M99").Value = -9718.666999999999
$ws4.Range("H105").Value = 5304.9
$ws4.Range("I105").Value = 2011.8
$ws4.Range("J105").Value = 8598
$ws4.Range("K105").Value = 2011.8
$ws4.Range("L105").Value = 8598
$ws4.Range("M105").Value = -264.8
$ws4.Range("N105").Value = -12092
$ws4.Range("H107").Value = 1645.4615
$ws4.Range("I107").Value = 1322.7059
$ws4.Range("J107").Value = 2255.111
$ws4.Range("K107").Value = 1322.7059
$ws4.Range("L107").Value = 2255.111
$ws4.Range("M107").Value = 597.2941000000001
$ws4.Range("N107").Value = -6095.111
$ws4.Range("H126").Value = 9941.666999999999
$ws4.Range("I126").Value = 11216.667
$ws4.Range("K126").Value = 33650.001
$ws4.Range("M126").Value = -31180.001

$ws5 = $wb.Worksheets.Item("CUL")
$ws5.Range("H25").Value = 1680.2
$ws5.Range("I25").Value = 1200
$ws5.Range("J25").Value = 2000.3334
$ws5.Range("K25").Value = 3600
$ws5.Range("L25").Value = 6001.0002
$ws5.Range("M25").Value = -3431
$ws5.Range("N25").Value = -6339.0002
$ws5.Range("H30").Value = 1680.2
$ws5.Range("I30").Value = 1200
$ws5.Range("J30").Value = 2000.3334
$ws5.Range("K30").Value = 3600
$ws5.Range("L30").Value = 6001.0002
$ws5.Range("M30").Value = -3498
$ws5.Range("N30").Value = -6205.0002
$ws5.Range("H75").Value = 3664
$ws5.Range("I75").Value = 2583
$ws5.Range("J75").Value = 4204.5
$ws5.Range("K75").Value = 7749
$ws5.Range("L75").Value = 12613.5
$ws5.Range("M75").Value = -6751
$ws5.Range("N75").Value = -14609.5
$ws5.Range("H78").Value = 3664
$ws5.Range("I78").Value = 2583
$ws5.Range("J78").Value = 4204.5
$ws5.Range("K78").Value = 23247
$ws5.Range("L78").Value = 37840.5
$ws5.Range("M78").Value = -18255
$ws5.Range("N78").Value = -47824.5
$ws5.Range("H81").Value = 5999.923
$ws5.Range("I81").Value = 2000
$ws5.Range("K81").Value = 6000
$ws5.Range("M81").Value = -4877
$ws5.Range("H84").Value = 5999.923
$ws5.Range("I84").Value = 2000
$ws5.Range("K84").Value = 18000
$ws5.Range("M84").Value = -12384
$ws5.Range("H87").Value = 3674.6667
$ws5.Range("J87").Value = 6996
$ws5.Range("L87").Value = 20988
$ws5.Range("N87").Value = -23484
$ws5.Range("H90").Value = 3674.6667
$ws5.Range("J90").Value = 6996
$ws5.Range("L90").Value = 62964
$ws5.Range("N90").Value = -75444
$ws5.Range("H103").Value = 185.7
$ws5.Range("I103").Value = 203.22223
$ws5.Range("J103").Value = 28
$ws5.Range("K103").Value = 609.66669
$ws5.Range("L103").Value = 84
$ws5.Range("M103").Value = 269.33331
$ws5.Range("N103").Value = -1842
$ws5.Range("H109").Value = 76390100
$ws5.Range("I109").Value = 83334440
$ws5.Range("K109").Value = 250003320
$ws5.Range("M109").Value = -250002280
$ws5.Range("H127").Value = 2000
$ws5.Range("J127").Value = 2000
$ws5.Range("L127").Value = 6000
$ws5.Range("N127").Value = -15920
$ws5.Range("H132").Value = 20777.555
$ws5.Range("I132").Value = 7500
$ws5.Range("J132").Value = 31399.6
$ws5.Range("K132").Value = 67500
$ws5.Range("L132").Value = 282596.4
$ws5.Range("M132").Value = -64970
$ws5.Range("N132").Value = -287656.4

$ws6 = $wb.Worksheets.Item("GSM")
$ws6.Range("H2").Value = 167.46153
$ws6.Range("I2").Value = 74.36364
$ws6.Range("K2").Value = 74.36364
$ws6.Range("M2").Value = 38.63636
$ws6.Range("H70").Value = 7022.4
$ws6.Range("I70").Value = 5999.56
$ws6.Range("K70").Value = 5999.56
$ws6.Range("M70").Value = -5729.56
$ws6.Range("H73").Value = 7022.4
$ws6.Range("I73").Value = 5999.56
$ws6.Range("K73").Value = 5999.56
$ws6.Range("M73").Value = -5063.56
$ws6.Range("H122").Value = 42685.926
$ws6.Range("I122").Value = 60360.445
$ws6.Range("J122").Value = 7336.8887
$ws6.Range("K122").Value = 181081.335
$ws6.Range("L122").Value = 22010.6661
$ws6.Range("M122").Value = -178631.335
$ws6.Range("N122").Value = -26910.6661
$ws6.Range("H132").Value = 3345.5
$ws6.Range("I132").Value = 3345.5
$ws6.Range("K132").Value = 10036.5
$ws6.Range("M132").Value = -7506.5

$ws7 = $wb.Worksheets.Item("LTW")
$ws7.Range("H22").Value = 2324.2
$ws7.Range("I22").Value = 2399.6667
$ws7.Range("J22").Value = 2211
$ws7.Range("K22").Value = 2399.6667
$ws7.Range("L22").Value = 2211
$ws7.Range("M22").Value = -2104.6667
$ws7.Range("N22").Value = -2801
$ws7.Range("H27").Value = 2324.2
$ws7.Range("I27").Value = 2399.6667
$ws7.Range("J27").Value = 2211
$ws7.Range("K27").Value = 2399.6667
$ws7.Range("L27").Value = 2211
$ws7.Range("M27").Value = -2292.6667
$ws7.Range("N27").Value = -2425
$ws7.Range("H132").Value = 12829247
$ws7.Range("I132").Value = 29417918
$ws7.Range("K132").Value = 88253754
$ws7.Range("M132").Value = -88251224

$ws8 = $wb.Worksheets.Item("WVR")
$ws8.Range("H4").Value = 3125.3845
$ws8.Range("I4").Value = 1851.5555
$ws8.Range("K4").Value = 1851.5555
$ws8.Range("M4").Value = -1738.5555
$ws8.Range("H5").Value = 7500000
$ws8.Range("I5").Value = 5000000
$ws8.Range("K5").Value = 5000000
$ws8.Range("M5").Value = -4999888
$ws8.Range("H81").Value = 8700350
$ws8.Range("I81").Value = 1792.579
$ws8.Range("J81").Value = 50018500
$ws8.Range("K81").Value = 3585.158
$ws8.Range("L81").Value = 100037000
$ws8.Range("M81").Value = -2524.158
$ws8.Range("N81").Value = -100039122
$ws8.Range("H84").Value = 8700350
$ws8.Range("I84").Value = 1792.579
$ws8.Range("J84").Value = 50018500
$ws8.Range("K84").Value = 17925.79
$ws8.Range("L84").Value = 500185000
$ws8.Range("M84").Value = -12621.79
$ws8.Range("N84").Value = -500195608
$ws8.Range("H132").Value = 17253174
$ws8.Range("I132").Value = 20835816
$ws8.Range("J132").Value = 56488.6
$ws8.Range("K132").Value = 62507448
$ws8.Range("L132").Value = 169465.8
$ws8.Range("M132").Value = -62504918
$ws8.Range("N132").Value = -174525.8
$ws8.Range("H135").Value = 72371.5
$ws8.Range("J135").Value = 72371.5
$ws8.Range("L135").Value = 72371.5
$ws8.Range("N135").Value = -82511.5
